$p = $ppt.ActivePresentation
Write-Output $p.Slides.Count
$s = $p.Slides.Item(1)
Write-Output $s.Shapes.Count
